# Add a new "Report Booking Statistic" API row to the api_config sheet,
# and rename the module/file-path style shared strings used across the
# whole sheet so they carry the "<module>/<name>" convention.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename existing module + excel_file values (A/G columns) -----------
$ws.Range("A2").Value = "user/functional"
$ws.Range("G2").Value = "user/Login.xlsx"

$ws.Range("A3").Value = "booking/functional"
$ws.Range("G3").Value = "booking/Booking_Quote_Fee.xlsx"

$ws.Range("A4").Value = "booking/functional"
$ws.Range("G4").Value = "booking/Create_Booking_Batch.xlsx"

$ws.Range("A5").Value = "booking/functional"
$ws.Range("G5").Value = "booking/Get_Booking_List_Select.xlsx"

$ws.Range("A6").Value = "booking/functional"

# --- new row data (row 7) ------------------------------------------------
$ws.Range("A7").Value = "booking/functional"
$ws.Range("B7").Value = "ReportBookingStatisticABC"
$ws.Range("C7").Value = "ReportBookingStatisticData"
$ws.Range("D7").Value = "/golf-cms/api/report/booking-statistic"
$ws.Range("E7").Value = "POST"
$ws.Range("F7").Value = "JSON"
$ws.Range("G7").Value = "booking/Report_Booking_Statistic.xlsx"
$ws.Range("H7").Value = "report_booking_statistic_template.json"
$ws.Range("I7").Value = "testReportBookingStatistic"

$ws.Rows.Item(7).RowHeight = 27

# --- column width tweaks (matches Excel's autofit after the new data) ----
$ws.Columns.Item(3).ColumnWidth = 26.7142857142857
$ws.Columns.Item(4).ColumnWidth = 35.4285714285714
$ws.Columns.Item(5).ColumnWidth = 16.2857142857143
$ws.Columns.Item(6).ColumnWidth = 19
$ws.Columns.Item(7).ColumnWidth = 36.4285714285714

# --- selection moves to H20, matching the saved view state ---------------
$ws.Range("H20").Select()
